$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 305.27274
$ws.Range("I2").Value = 177.2
$ws.Range("J2").Value = 412
$ws.Range("K2").Value = 177.2
$ws.Range("L2").Value = 412
$ws.Range("M2").Value = -64.19999999999999
$ws.Range("N2").Value = -638

# Row 45 (Leve Item ID 4585)
$ws.Range("H45").Value = 375
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 289.625
$ws.Range("I80").Value = 358.57144
$ws.Range("J80").Value = 236
$ws.Range("K80").Value = 1075.71432
$ws.Range("L80").Value = 708
$ws.Range("M80").Value = -77.71432000000004
$ws.Range("N80").Value = -2704

# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 289.625
$ws.Range("I83").Value = 358.57144
$ws.Range("J83").Value = 236
$ws.Range("K83").Value = 3227.14296
$ws.Range("L83").Value = 2124
$ws.Range("M83").Value = 1764.85704
$ws.Range("N83").Value = -12108

# Row 105 (Leve Item ID 18668)
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 965
$ws.Range("I127").Value = 982.1429000000001
$ws.Range("J127").Value = 845
$ws.Range("K127").Value = 2946.4287
$ws.Range("L127").Value = 2535
$ws.Range("M127").Value = 2013.5713
$ws.Range("N127").Value = -12455

$ws = $wb.Worksheets.Item("ARM")
# Row 13 (Leve Item ID 2656)
$ws.Range("H13").Value = 1135.8572
$ws.Range("J13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("N13").Value = -2788

# Row 14 (Leve Item ID 2673)
$ws.Range("H14").Value = 3866.6667
$ws.Range("I14").Value = 100
$ws.Range("J14").Value = 5750
$ws.Range("K14").Value = 100
$ws.Range("L14").Value = 5750
$ws.Range("M14").Value = 75
$ws.Range("N14").Value = -6100

# Row 19 (Leve Item ID 3550)
$ws.Range("H19").Value = 800
$ws.Range("I19").Value = 800
$ws.Range("K19").Value = 800
$ws.Range("M19").Value = -571

# Row 22 (Leve Item ID 2497)
$ws.Range("H22").Value = 4140
$ws.Range("I22").Value = 2675
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 2675
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -2376
$ws.Range("N22").Value = -10598

# Row 25 (Leve Item ID 2471)
$ws.Range("H25").Value = 1750
$ws.Range("I25").Value = 1750
$ws.Range("K25").Value = 1750
$ws.Range("M25").Value = -1348

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 5054.6
$ws.Range("I45").Value = 2256.1667
$ws.Range("K45").Value = 2256.1667
$ws.Range("M45").Value = -1879.1667

$ws = $wb.Worksheets.Item("CRP")
# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 10110.75
$ws.Range("I86").Value = 10147.667
$ws.Range("K86").Value = 10147.667
$ws.Range("M86").Value = -9024.666999999999

# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 10110.75
$ws.Range("I89").Value = 10147.667
$ws.Range("K89").Value = 50738.335
$ws.Range("M89").Value = -45122.335

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2367.5557
$ws.Range("I132").Value = 2367.5557
$ws.Range("K132").Value = 7102.6671
$ws.Range("M132").Value = -4572.6671

$ws = $wb.Worksheets.Item("CUL")
# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 68827.27
$ws.Range("I55").Value = 550
$ws.Range("K55").Value = 1650
$ws.Range("M55").Value = -1473

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 27781428
$ws.Range("I131").Value = 83334320
$ws.Range("K131").Value = 250002960
$ws.Range("M131").Value = -249997920

# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 1024.3
$ws.Range("I137").Value = 1067.7894
$ws.Range("J137").Value = 198
$ws.Range("K137").Value = 3203.3682
$ws.Range("L137").Value = 594
$ws.Range("M137").Value = 1896.6318
$ws.Range("N137").Value = -10794

$ws = $wb.Worksheets.Item("GSM")
# Row 62 (Leve Item ID 11983)
$ws.Range("H62").Value = 44000
$ws.Range("I62").Value = 44000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 44000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -43314
$ws.Range("N62").ClearContents()

# Row 65 (Leve Item ID 11983)
$ws.Range("H65").Value = 44000
$ws.Range("I65").Value = 44000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 132000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -128568
$ws.Range("N65").ClearContents()

# Row 134 (Leve Item ID 42064)
$ws.Range("H134").Value = 52717.668
$ws.Range("J134").Value = 52717.668
$ws.Range("L134").Value = 158153.004
$ws.Range("N134").Value = -163223.004

$ws = $wb.Worksheets.Item("LTW")
# Row 20 (Leve Item ID 4308)
$ws.Range("H20").Value = 52499.9
$ws.Range("J20").Value = 52499.9
$ws.Range("L20").Value = 52499.9
$ws.Range("N20").Value = -52951.9

# Row 63 (Leve Item ID 12006)
$ws.Range("H63").Value = 49384.715
$ws.Range("J63").Value = 49948.832
$ws.Range("L63").Value = 49948.832
$ws.Range("N63").Value = -51446.832

# Row 66 (Leve Item ID 12006)
$ws.Range("H66").Value = 49384.715
$ws.Range("J66").Value = 49948.832
$ws.Range("L66").Value = 149846.496
$ws.Range("N66").Value = -157334.496

# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 322.6
$ws.Range("I82").Value = 178.5
$ws.Range("K82").Value = 178.5
$ws.Range("M82").Value = 182.5

# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 322.6
$ws.Range("I85").Value = 178.5
$ws.Range("K85").Value = 178.5
$ws.Range("M85").Value = 1069.5

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 4379.4346
$ws.Range("I136").Value = 3633
$ws.Range("K136").Value = 10899
$ws.Range("M136").Value = -8349

$ws = $wb.Worksheets.Item("WVR")
# Row 30 (Leve Item ID 2700)
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# Row 46 (Leve Item ID 42037)
$ws.Range("H46").Value = 69598.8
$ws.Range("I46").Value = 50000
$ws.Range("K46").Value = 50000
$ws.Range("M46").Value = -49769

# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 1047.7693
$ws.Range("I100").Value = 984.7273
$ws.Range("K100").Value = 1969.4546
$ws.Range("M100").Value = -1428.4546

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 1697.55
$ws.Range("I126").Value = 1697.55
$ws.Range("K126").Value = 5092.65
$ws.Range("M126").Value = -2622.65

# Row 134 (Leve Item ID 42037)
$ws.Range("H134").Value = 69598.8
$ws.Range("I134").Value = 50000
$ws.Range("K134").Value = 150000
$ws.Range("M134").Value = -147465

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 7003.0527
$ws.Range("I136").Value = 6002.75
$ws.Range("K136").Value = 18008.25
$ws.Range("M136").Value = -15458.25
